$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price), E (Volume 1h) and G (Hora) hold numeric-looking text
# (e.g. "292.08", "-0.56%", "19") in the source workbook - the cells are
# stored as text, not numbers. A plain `.Value = "19"` assignment lets
# Excel auto-convert such digit-only strings to real numbers, so the cell
# is switched to Text format first; this keeps the stored cell as a text
# string, matching the workbook being edited.
function Set-TextValue($sheet, $addr, $val) {
    $cell = $sheet.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
}

Set-TextValue $ws 'D2' '292.08'
Set-TextValue $ws 'E2' '-0.56%'
Set-TextValue $ws 'G2' '19'
Set-TextValue $ws 'D3' '30.88'
Set-TextValue $ws 'E3' '-0.56%'
Set-TextValue $ws 'G3' '19'
Set-TextValue $ws 'D4' '4.891'
Set-TextValue $ws 'E4' '0.12%'
Set-TextValue $ws 'G4' '19'
Set-TextValue $ws 'D5' '0.07277'
Set-TextValue $ws 'G5' '19'
Set-TextValue $ws 'D6' '2.284'
Set-TextValue $ws 'E6' '24.95%'
Set-TextValue $ws 'G6' '19'
Set-TextValue $ws 'D7' '7.681'
Set-TextValue $ws 'E7' '0.06%'
Set-TextValue $ws 'G7' '19'
Set-TextValue $ws 'D8' '3.717'
Set-TextValue $ws 'E8' '-1.33%'
Set-TextValue $ws 'G8' '19'
Set-TextValue $ws 'D9' '0.8990'
Set-TextValue $ws 'E9' '-0.87%'
Set-TextValue $ws 'G9' '19'
Set-TextValue $ws 'D10' '0.1668'
Set-TextValue $ws 'E10' '0.85%'
Set-TextValue $ws 'G10' '19'
Set-TextValue $ws 'D11' '0.07912'
Set-TextValue $ws 'E11' '4.46%'
Set-TextValue $ws 'G11' '19'
Set-TextValue $ws 'D12' '0.08021'
Set-TextValue $ws 'E12' '-1.85%'
Set-TextValue $ws 'G12' '19'
Set-TextValue $ws 'D13' '0.03098'
Set-TextValue $ws 'E13' '3.39%'
Set-TextValue $ws 'G13' '19'
Set-TextValue $ws 'E14' '0.34%'
Set-TextValue $ws 'G14' '19'
Set-TextValue $ws 'D15' '0.001497'
Set-TextValue $ws 'E15' '-0.46%'
Set-TextValue $ws 'G15' '19'
Set-TextValue $ws 'D16' '0.005774'
Set-TextValue $ws 'E16' '1.98%'
Set-TextValue $ws 'G16' '19'
Set-TextValue $ws 'D17' '3.482'
Set-TextValue $ws 'E17' '0.65%'
Set-TextValue $ws 'G17' '19'
Set-TextValue $ws 'E18' '-0.89%'
Set-TextValue $ws 'G18' '19'
Set-TextValue $ws 'D19' '0.3323'
Set-TextValue $ws 'E19' '1.54%'
Set-TextValue $ws 'G19' '19'
Set-TextValue $ws 'D20' '0.1299'
Set-TextValue $ws 'E20' '-0.58%'
Set-TextValue $ws 'G20' '19'
Set-TextValue $ws 'D21' '4.031'
Set-TextValue $ws 'E21' '-7.74%'
Set-TextValue $ws 'G21' '19'
Set-TextValue $ws 'D22' '0.2098'
Set-TextValue $ws 'E22' '4.76%'
Set-TextValue $ws 'G22' '19'
Set-TextValue $ws 'D23' '0.04510'
Set-TextValue $ws 'E23' '0.61%'
Set-TextValue $ws 'G23' '19'
Set-TextValue $ws 'E24' '-1.44%'
Set-TextValue $ws 'G24' '19'
Set-TextValue $ws 'D25' '0.004657'
Set-TextValue $ws 'E25' '15.34%'
Set-TextValue $ws 'G25' '19'
Set-TextValue $ws 'E26' '3.80%'
Set-TextValue $ws 'G26' '19'
Set-TextValue $ws 'D27' '0.0003390'
Set-TextValue $ws 'G27' '19'
Set-TextValue $ws 'G28' '19'
Set-TextValue $ws 'G29' '19'
Set-TextValue $ws 'G30' '19'
Set-TextValue $ws 'G31' '19'
Set-TextValue $ws 'G32' '19'
Set-TextValue $ws 'G33' '19'
Set-TextValue $ws 'G34' '19'
Set-TextValue $ws 'G35' '19'
Set-TextValue $ws 'G36' '19'
Set-TextValue $ws 'G37' '19'
Set-TextValue $ws 'G38' '19'
Set-TextValue $ws 'D39' '0.01583'
Set-TextValue $ws 'E39' '-3.99%'
Set-TextValue $ws 'G39' '19'
Set-TextValue $ws 'D40' '0.04382'
Set-TextValue $ws 'E40' '-0.32%'
Set-TextValue $ws 'G40' '19'
Set-TextValue $ws 'D41' '0.007304'
Set-TextValue $ws 'E41' '-1.83%'
Set-TextValue $ws 'G41' '19'
Set-TextValue $ws 'D42' '0.009755'
Set-TextValue $ws 'G42' '19'
Set-TextValue $ws 'D43' '0.1313'
Set-TextValue $ws 'E43' '-0.48%'
Set-TextValue $ws 'G43' '19'
Set-TextValue $ws 'D44' '0.002060'
Set-TextValue $ws 'E44' '-2.56%'
Set-TextValue $ws 'G44' '19'
Set-TextValue $ws 'D45' '0.009360'
Set-TextValue $ws 'E45' '-16.29%'
Set-TextValue $ws 'G45' '19'
Set-TextValue $ws 'E46' '-3.76%'
Set-TextValue $ws 'G46' '19'
Set-TextValue $ws 'D47' '0.00000000750'
Set-TextValue $ws 'E47' '-0.19%'
Set-TextValue $ws 'G47' '19'
$ws.Range('B48').Value = 'BOLO'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ogrGe0dEab+bolo-bolo'
Set-TextValue $ws 'D48' '2.241'
Set-TextValue $ws 'E48' '4.56%'
Set-TextValue $ws 'G48' '19'
$ws.Range('B49').Value = 'CoinbaseStockToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin'
Set-TextValue $ws 'D49' '0.002895'
Set-TextValue $ws 'E49' '20.43%'
Set-TextValue $ws 'G49' '19'
Set-TextValue $ws 'D50' '0.00002099'
Set-TextValue $ws 'E50' '-0.19%'
Set-TextValue $ws 'G50' '19'
Set-TextValue $ws 'D51' '0.0001999'
Set-TextValue $ws 'E51' '-0.19%'
Set-TextValue $ws 'G51' '19'
